# Insert a new weekly price record above row 112 (Terminal Hortofrutícola
# Agro Chillán - Alcachofa), pushing the existing rows 112-123 down to
# 113-124 and extending the sheet dimension to A1:R124.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 112..123 down one position, creating a fresh blank row 112.
$ws.Rows.Item(112).Insert()

# Populate the new row 112 with the latest weekly data point.
$ws.Cells.Item(112, 1).Value = 7
$ws.Cells.Item(112, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(112, 3).Value = "Ñuble"
$ws.Cells.Item(112, 4).Value = 45166
$ws.Cells.Item(112, 5).Value = 16
$ws.Cells.Item(112, 6).Value = 100112013
$ws.Cells.Item(112, 7).Value = "Alcachofa"
$ws.Cells.Item(112, 8).Value = "Madrigal"
$ws.Cells.Item(112, 9).Value = "Primera"
$ws.Cells.Item(112, 10).Value = 100
$ws.Cells.Item(112, 11).Value = 12000
$ws.Cells.Item(112, 12).Value = 12000
$ws.Cells.Item(112, 13).Value = 12000
$ws.Cells.Item(112, 14).Value = "`$/caja 40 unidades"
$ws.Cells.Item(112, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(112, 16).Value = 300
$ws.Cells.Item(112, 17).Value = 40
$ws.Cells.Item(112, 18).Value = "Hortaliza"
